$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title paragraph: "Project Title (TBD)" -> underlined new title
# ------------------------------------------------------------------
$d.Content.Find.Execute("Project Title (TBD)", $true, $false, $false, $false, $false, $true, 1, $false, "How Education Impacts Income Across the United States", 2)
$pTitle = $d.Paragraphs.Item(1)
$pTitle.Range.Font.Underline = 1

# ------------------------------------------------------------------
# 2) The blank paragraph right after the title becomes centered +
#    underlined (still empty of visible text).
# ------------------------------------------------------------------
$pBlank = $d.Paragraphs.Item(2)
$pBlank.Alignment = 1
$blankRange = $pBlank.Range
$blankRange.InsertBefore("X")
$pBlank2 = $d.Paragraphs.Item(2)
$pBlank2.Range.Font.Underline = 1
$tmpRange = $d.Range($pBlank2.Range.Start, $pBlank2.Range.Start + 1)
$tmpRange.Text = ""

# ------------------------------------------------------------------
# 3) Research questions: swap question 1 & 2 (each now ends with a
#    separate "?" run), and replace question 3 entirely. Done before
#    the "Lexie Fallow" insertion below so the paragraph indices
#    (13/14/15) are still the original ones.
# ------------------------------------------------------------------
$q1 = $d.Paragraphs.Item(13)
$q1.Range.Text = "What is the correlation between education levels (high education/no higher education) and income"
$q1b = $d.Paragraphs.Item(13)
$q1b.Range.InsertAfter("?")
$q1c = $d.Paragraphs.Item(13)
$q1mark = $d.Range($q1c.Range.End - 2, $q1c.Range.End - 1)
$q1mark.Bold = 1
$q1mark.Bold = 0

$q2 = $d.Paragraphs.Item(14)
$q2.Range.Text = "What degree type (associates or bachelors) has the greatest impact on income"
$q2b = $d.Paragraphs.Item(14)
$q2b.Range.InsertAfter("?")
$q2c = $d.Paragraphs.Item(14)
$q2mark = $d.Range($q2c.Range.End - 2, $q2c.Range.End - 1)
$q2mark.Bold = 1
$q2mark.Bold = 0

$q3 = $d.Paragraphs.Item(15)
$q3.Range.Text = "Which US region has the highest amount of people with higher education degrees and does this region also have the highest income?"

# ------------------------------------------------------------------
# 4) Team members: insert "Lexie Fallow" as a new list item right
#    before the existing blank ListParagraph line.
# ------------------------------------------------------------------
$pKudirat = $d.Paragraphs.Item(7)
$pKudirat.Range.InsertParagraphAfter()
$pLexie = $d.Paragraphs.Item(8)
$pLexie.Range.InsertBefore("Lexie Fallow")

Write-Host "done"
